# 2-2-2.xlsx "Add files via upload" edit
#
# The organization's contact-details block (rows 6-10 in column B) was
# refreshed with new staff/contact information, and the previously
# selected cell moved from B25 to B8.
#
# Contact info is entered in this particular order so that the workbook's
# shared-string table ends up built in the same sequence as the reference
# edit (Калымбетова Ы.И., phone, organization, website, e-mail) even
# though the cells themselves are written in sheet order afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value  = "Калымбетова Ы.И."
$ws.Range("B9").Value  = "(0312) 32 46 55"
$ws.Range("B6").Value  = "Национальный статистический комитет КР (Управление статистики домашних хозяйств) в рамках глобальной программы MICS ЮНИСЕФ"
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Range("B8").Value  = "yryskan.kalymbetova@gmail.com"

# Retyped cells pick up a (freshly materialised) Cyrillic-capable Calibri
# font; the e-mail cell (B8) was left with its original formatting.
$ws.Range("B7").Font.Name  = "Calibri"
$ws.Range("B9").Font.Name  = "Calibri"
$ws.Range("B10").Font.Name = "Calibri"
$ws.Range("B6").Font.Name  = "Calibri"

# Selection moved to the e-mail cell.
$ws.Range("B8").Select()
